$p = $ppt.ActivePresentation

# --- Slide 1: split subtitle run into "Data " + "Science vs Data Analytics" ---
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2).TextFrame.TextRange
$subtitle.Text = "Data Science vs Data Analytics"
$firstWord = $subtitle.Characters(1, 5)
$firstWord.Text = "Data "

# --- Slide 2: fill the empty content placeholder with the bullet list ---
$s2 = $p.Slides.Item(2)
$body = $s2.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Asking questions (formulating hypothesis), answers to which solve known problems or unearth unknown solutions that in turn drive business value,"
$body.InsertAfter("`rDefining the data needed or working with an existing data set and employing tools (computer science based) to collect, store and explore such data generally in huge volume & variety (often more than 1 TB and 1000s of dimensions),")
$body.InsertAfter("`rIdentifying the type of analysis to be done to get to the answers and performing such analysis by implementing various algorithms/tools (statistics based), often in a distributed and parallel architecture,")
$body.InsertAfter("`rCommunicating the insights gathered from the analysis in the form of simple stories/visualizations/dashboards (the Data Product) that a non-data scientist can understand and build conversation out of it. (It should be kept in mind that a product can also be an piece of code that is internal to a company and is used by various departments. The presentation, maintenance, scalability, etc of the code are then the product features, which is often not practiced in many organizations)")
$body.InsertAfter("`rBuilding a higher level abstraction that does steps 2-3-4 in an autonomous way, analyzing & taking actions on new data as they are fed to the system.")

# Turn on "shrink text on overflow" autofit for the now much-longer body text.
$s2.Shapes.Item(2).TextFrame.AutoSize = 2
